$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add column R (year 2021) to the header row (row 4) ---
# Copy formatting from the existing Q4 cell (year 2020 header) so the new
# cell reuses the same style (numFmtId 0 / fontId 20 / borderId 1 / xfId 1,
# horizontal+vertical center) instead of creating a brand-new one.
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial(-4122)
$ws.Range("R4").Value = 2021

# --- Add column R (growth rate value) to the data row (row 5) ---
# Start from Q5's formatting (Times New Roman 9, border, right/vertical
# center) and then apply the new "0.0" number format, which causes the
# interop runtime to create the new numFmt (167) and cellXfs entry used by
# this single cell.
$ws.Range("Q5").Copy()
$ws.Range("R5").PasteSpecial(-4122)
$ws.Range("R5").Value = 102.20441221981518
$ws.Range("R5").NumberFormat = "0.0"

# Clear clipboard marching ants / copy mode
$excel.CutCopyMode = 0

# --- Update the active selection shown in the sheet view ---
[void]$ws.Range("S9").Select()
